$p = $ppt.ActivePresentation

# --- Slide 1 ("Semana 2" -> "Semana " + "1" typed as two runs) ---
$s1 = $p.Slides.Item(1)
$title1 = $s1.Shapes.Item(1).TextFrame.TextRange
$title1.Text = "Semana "
[void]$title1.InsertAfter("1")

# --- Slide 19 (re-join the accidentally split "subred" / "?" runs) ---
$s19 = $p.Slides.Item(19)
$title19 = $s19.Shapes.Item(1).TextFrame.TextRange
$title19.Text = "TEMP"
$title19.Text = "¿Cuál es la dirección de subred?"

# --- Slide 21 (re-join the accidentally split "El " / "hilo tiene..." runs) ---
$s21 = $p.Slides.Item(21)
$body21 = $s21.Shapes.Item(2).TextFrame.TextRange
$para5 = $body21.Paragraphs(5, 1)
$para5.Text = "TEMP"
$para5b = $body21.Paragraphs(5, 1)
$para5b.Text = "El hilo tiene que reportar a la interfaz cada vez que encuentre un host"
